$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1 ("Save"), matching the style of the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new "Save" column for rows 2-9
$saveValues = @(0, 1, 1, 1, 1, 0, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
